$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02572360017603
$ws.Range("D2").Value = 1.036581664324456
$ws.Range("E2").Value = 1.026019270175057
$ws.Range("F2").Value = 1.02422791284009
$ws.Range("I2").Value = 1.035804277293185
$ws.Range("J2").Value = 1.030891027493712
$ws.Range("K2").Value = 1.039374995436167
$ws.Range("L2").Value = 1.028843141018787
$ws.Range("M2").Value = 1.027057031097957
$ws.Range("N2").Value = 1.014245111639921
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026975438004378
$ws.Range("D3").Value = 1.037252778767501
$ws.Range("E3").Value = 1.027090843599717
$ws.Range("F3").Value = 1.026119928692946
$ws.Range("I3").Value = 1.036122666004603
$ws.Range("J3").Value = 1.031780762318228
$ws.Range("K3").Value = 1.03985610845312
$ws.Range("L3").Value = 1.029721424141658
$ws.Range("M3").Value = 1.028753143468405
$ws.Range("N3").Value = 1.014546881665385
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02778410373823
$ws.Range("D4").Value = 1.037686136142086
$ws.Range("E4").Value = 1.027783374195255
$ws.Range("F4").Value = 1.027342507476574
$ws.Range("I4").Value = 1.036326640842956
$ws.Range("J4").Value = 1.032354673639352
$ws.Range("K4").Value = 1.04016584423969
$ws.Range("L4").Value = 1.03028829704306
$ws.Range("M4").Value = 1.029848565064724
$ws.Range("N4").Value = 1.014741354752175
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028123747221922
$ws.Range("D5").Value = 1.037868105041781
$ws.Range("E5").Value = 1.028074314502005
$ws.Range("F5").Value = 1.027856090312354
$ws.Range("I5").Value = 1.036411903717188
$ws.Range("J5").Value = 1.032595517368804
$ws.Range("K5").Value = 1.040295680627231
$ws.Range("L5").Value = 1.030526270392548
$ws.Range("M5").Value = 1.030308595141527
$ws.Range("N5").Value = 1.014822922784961
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028180756269361
$ws.Range("D6").Value = 1.037898645856174
$ws.Range("E6").Value = 1.028123153064933
$ws.Range("F6").Value = 1.027942300689823
$ws.Range("I6").Value = 1.036426191137607
$ws.Range("J6").Value = 1.032635931071147
$ws.Range("K6").Value = 1.040317458672908
$ws.Range("L6").Value = 1.030566207328099
$ws.Range("M6").Value = 1.030385808111558
$ws.Range("N6").Value = 1.014836607401
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027788643320338
$ws.Range("D7").Value = 1.037688568461772
$ws.Range("E7").Value = 1.027787262533429
$ws.Range("F7").Value = 1.027349371508698
$ws.Range("I7").Value = 1.036327782046318
$ws.Range("J7").Value = 1.032357893484316
$ws.Range("K7").Value = 1.040167580598341
$ws.Range("L7").Value = 1.030291478186317
$ws.Range("M7").Value = 1.029854713900763
$ws.Range("N7").Value = 1.014742445407209
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026146948812308
$ws.Range("D8").Value = 1.036808656995076
$ws.Range("E8").Value = 1.026381591479339
$ws.Range("F8").Value = 1.024867684805415
$ws.Range("I8").Value = 1.035912302175802
$ws.Range("J8").Value = 1.031192093979305
$ws.Range("K8").Value = 1.039537916796022
$ws.Range("L8").Value = 1.0291402601713
$ws.Range("M8").Value = 1.027630677395117
$ws.Range("N8").Value = 1.014347261024156
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023243441563972
$ws.Range("D9").Value = 1.035251231467727
$ws.Range("E9").Value = 1.023897948968765
$ws.Range("F9").Value = 1.020481151839236
$ws.Range("I9").Value = 1.035164465641507
$ws.Range("J9").Value = 1.029123793099687
$ws.Range("K9").Value = 1.038416253760886
$ws.Range("L9").Value = 1.027100517269082
$ws.Range("M9").Value = 1.023695199293853
$ws.Range("N9").Value = 1.013644772216539
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02130028550333
$ws.Range("D10").Value = 1.034208260329266
$ws.Range("E10").Value = 1.022237464803667
$ws.Range("F10").Value = 1.017546886943639
$ws.Range("I10").Value = 1.034655269503343
$ws.Range("J10").Value = 1.02773525415324
$ws.Range("K10").Value = 1.037660272816626
$ws.Range("L10").Value = 1.025732966870358
$ws.Range("M10").Value = 1.021059738379953
$ws.Range("N10").Value = 1.013172250755952
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020457026488083
$ws.Range("D11").Value = 1.033755519915431
$ws.Range("E11").Value = 1.021517282901658
$ws.Range("F11").Value = 1.016273779085741
$ws.Range("I11").Value = 1.034432240761018
$ws.Range("J11").Value = 1.027131650972387
$ws.Range("K11").Value = 1.0373309641242
$ws.Range("L11").Value = 1.02513891995874
$ws.Range("M11").Value = 1.019915585300954
$ws.Range("N11").Value = 1.012966630792475
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020143516527129
$ws.Range("D12").Value = 1.033587181607622
$ws.Range("E12").Value = 1.021249593062902
$ws.Range("F12").Value = 1.015800490477908
$ws.Range("I12").Value = 1.03434901420156
$ws.Range("J12").Value = 1.026907086846621
$ws.Range("K12").Value = 1.037208347808953
$ws.Range("L12").Value = 1.024917976262301
$ws.Range("M12").Value = 1.019490133438023
$ws.Range("N12").Value = 1.012890100305962
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020210778527949
$ws.Range("D13").Value = 1.033623298472199
$ws.Range("E13").Value = 1.021307021764182
$ws.Range("F13").Value = 1.015902030839628
$ws.Range("I13").Value = 1.034366883973172
$ws.Range("J13").Value = 1.02695527294103
$ws.Range("K13").Value = 1.037234662859573
$ws.Range("L13").Value = 1.024965382550475
$ws.Range("M13").Value = 1.01958141547005
$ws.Range("N13").Value = 1.012906523352904
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020431117530135
$ws.Range("D14").Value = 1.033741608499793
$ws.Range("E14").Value = 1.021495159320778
$ws.Range("F14").Value = 1.016234665149557
$ws.Range("I14").Value = 1.034425369065628
$ws.Range("J14").Value = 1.027113095784705
$ws.Range("K14").Value = 1.037320834675757
$ws.Range("L14").Value = 1.02512066258967
$ws.Range("M14").Value = 1.019880426836281
$ws.Range("N14").Value = 1.012960307910888
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020566837562112
$ws.Range("D15").Value = 1.033814480606632
$ws.Range("E15").Value = 1.021611052790985
$ws.Range("F15").Value = 1.016439558677023
$ws.Range("I15").Value = 1.034461352761047
$ws.Range("J15").Value = 1.027210287923165
$ws.Range("K15").Value = 1.037373888658639
$ws.Range("L15").Value = 1.025216297419224
$ws.Range("M15").Value = 1.020064595887679
$ws.Range("N15").Value = 1.012993425888041
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021356210082457
$ws.Range("D16").Value = 1.034238283386516
$ws.Range("E16").Value = 1.022285235612798
$ws.Range("F16").Value = 1.017631323769244
$ws.Range("I16").Value = 1.034670017460285
$ws.Range("J16").Value = 1.027775263233544
$ws.Range("K16").Value = 1.037682086457771
$ws.Range("L16").Value = 1.025772351635543
$ws.Range("M16").Value = 1.021135608000976
$ws.Range("N16").Value = 1.013185875555553
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021850860075368
$ws.Range("D17").Value = 1.034503821213834
$ws.Range("E17").Value = 1.022707813088559
$ws.Range("F17").Value = 1.018378191887179
$ws.Range("I17").Value = 1.034800225215571
$ws.Range("J17").Value = 1.02812902288442
$ws.Range("K17").Value = 1.037874884048934
$ws.Range("L17").Value = 1.026120640978285
$ws.Range("M17").Value = 1.021806617219224
$ws.Range("N17").Value = 1.013306321281206
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022139202073031
$ws.Range("D18").Value = 1.034658596266812
$ws.Range("E18").Value = 1.022954181725723
$ws.Range("F18").Value = 1.018813582320847
$ws.Range("I18").Value = 1.03487592790163
$ws.Range("J18").Value = 1.02833513785962
$ws.Range("K18").Value = 1.037987150221741
$ws.Range("L18").Value = 1.02632361034649
$ws.Range("M18").Value = 1.022197718901869
$ws.Range("N18").Value = 1.01337647740497
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022237489084969
$ws.Range("D19").Value = 1.034711352158448
$ws.Range("E19").Value = 1.023038167928619
$ws.Range("F19").Value = 1.018961998123803
$ws.Range("I19").Value = 1.03490169900373
$ws.Range("J19").Value = 1.028405379378325
$ws.Range("K19").Value = 1.038025397999778
$ws.Range("L19").Value = 1.026392786935069
$ws.Range("M19").Value = 1.02233102628846
$ws.Range("N19").Value = 1.013400382283802
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02179780737414
$ws.Range("D20").Value = 1.03447534276632
$ws.Range("E20").Value = 1.022662486316031
$ws.Range("F20").Value = 1.01829808547923
$ws.Range("I20").Value = 1.034786280540006
$ws.Range("J20").Value = 1.028091091343211
$ws.Range("K20").Value = 1.037854218292892
$ws.Range("L20").Value = 1.026083291677687
$ws.Range("M20").Value = 1.021734654037275
$ws.Range("N20").Value = 1.013293408725592
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020366241134808
$ws.Range("D21").Value = 1.033706773856713
$ws.Range("E21").Value = 1.021439762567866
$ws.Range("F21").Value = 1.016136723800207
$ws.Range("I21").Value = 1.034408157270296
$ws.Range("J21").Value = 1.027066630841987
$ws.Range("K21").Value = 1.037295467416797
$ws.Range("L21").Value = 1.025074944481038
$ws.Range("M21").Value = 1.019792388325737
$ws.Range("N21").Value = 1.012944473965129
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019464500489519
$ws.Range("D22").Value = 1.033222557893122
$ws.Range("E22").Value = 1.02066993099968
$ws.Range("F22").Value = 1.014775469419329
$ws.Range("I22").Value = 1.034168195055629
$ws.Range("J22").Value = 1.026420432135478
$ws.Range("K22").Value = 1.036942442974824
$ws.Range("L22").Value = 1.024439286748279
$ws.Range("M22").Value = 1.018568525776715
$ws.Range("N22").Value = 1.01272419259779
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019942689667876
$ws.Range("D23").Value = 1.033479343769951
$ws.Range("E23").Value = 1.02107813502901
$ws.Range("F23").Value = 1.015497321555955
$ws.Range("I23").Value = 1.034295614650925
$ws.Range("J23").Value = 1.026763193102245
$ws.Range("K23").Value = 1.037129750984073
$ws.Range("L23").Value = 1.024776420730132
$ws.Range("M23").Value = 1.019217577693357
$ws.Range("N23").Value = 1.012841053022618
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021821780158686
$ws.Range("D24").Value = 1.034488211286717
$ws.Range("E24").Value = 1.022682967884314
$ws.Range("F24").Value = 1.01833428287253
$ws.Range("I24").Value = 1.034792582296473
$ws.Range("J24").Value = 1.028108231675073
$ws.Range("K24").Value = 1.037863556843681
$ws.Range("L24").Value = 1.026100168780686
$ws.Range("M24").Value = 1.021767171985789
$ws.Range("N24").Value = 1.013299243655824
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023995362144953
$ws.Range("D25").Value = 1.03565468629461
$ws.Range("E25").Value = 1.024540844887865
$ws.Range("F25").Value = 1.021616860100078
$ws.Range("I25").Value = 1.035359668663282
$ws.Range("J25").Value = 1.029660184669696
$ws.Range("K25").Value = 1.038707672327926
$ws.Range("L25").Value = 1.027629184054075
$ws.Range("M25").Value = 1.024714640409229
$ws.Range("N25").Value = 1.01382711610304
